# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# updates (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR Leve-profit worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (ALC)
$ws.Cells.Item(32, 8).Value = 820.2
$ws.Cells.Item(32, 9).Value = 700
$ws.Cells.Item(32, 10).Value = 850.25
$ws.Cells.Item(32, 11).Value = 700
$ws.Cells.Item(32, 12).Value = 850.25
$ws.Cells.Item(32, 13).Value = -374
$ws.Cells.Item(32, 14).Value = -1502.25

# Row 52 (ALC)
$ws.Cells.Item(52, 8).Value = 299797
$ws.Cells.Item(52, 10).Value = 201531.33
$ws.Cells.Item(52, 12).Value = 604593.99
$ws.Cells.Item(52, 14).Value = -604913.99

# Row 86 (ALC)
$ws.Cells.Item(86, 8).Value = 1725.875
$ws.Cells.Item(86, 9).Value = 1686.1428
$ws.Cells.Item(86, 10).Value = 2004
$ws.Cells.Item(86, 11).Value = 1686.1428
$ws.Cells.Item(86, 12).Value = 2004
$ws.Cells.Item(86, 13).Value = -563.1428000000001
$ws.Cells.Item(86, 14).Value = -4250

# Row 89 (ALC)
$ws.Cells.Item(89, 8).Value = 1725.875
$ws.Cells.Item(89, 9).Value = 1686.1428
$ws.Cells.Item(89, 10).Value = 2004
$ws.Cells.Item(89, 11).Value = 8430.714
$ws.Cells.Item(89, 12).Value = 10020
$ws.Cells.Item(89, 13).Value = -2814.714
$ws.Cells.Item(89, 14).Value = -21252

# Row 127 (ALC)
$ws.Cells.Item(127, 8).Value = 1488.5555
$ws.Cells.Item(127, 9).Value = 674.25
$ws.Cells.Item(127, 10).Value = 2140
$ws.Cells.Item(127, 11).Value = 2022.75
$ws.Cells.Item(127, 12).Value = 6420
$ws.Cells.Item(127, 13).Value = 2937.25
$ws.Cells.Item(127, 14).Value = -16340

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 2501586.2
$ws.Cells.Item(137, 9).Value = 4349079.5
$ws.Cells.Item(137, 10).Value = 2036.6471
$ws.Cells.Item(137, 11).Value = 13047238.5
$ws.Cells.Item(137, 12).Value = 6109.9413
$ws.Cells.Item(137, 13).Value = -13044688.5
$ws.Cells.Item(137, 14).Value = -11209.9413

# Row 138 (ALC)
$ws.Cells.Item(138, 8).Value = 3335360.2
$ws.Cells.Item(138, 9).Value = 1590.5862
$ws.Cells.Item(138, 10).Value = 7939137.5
$ws.Cells.Item(138, 11).Value = 4771.7586
$ws.Cells.Item(138, 12).Value = 23817412.5
$ws.Cells.Item(138, 13).Value = 368.2413999999999
$ws.Cells.Item(138, 14).Value = -23827692.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Cells.Item(32, 8).Value = 814.71
$ws.Cells.Item(32, 9).Value = 710.55054
$ws.Cells.Item(32, 10).Value = 1657.4546
$ws.Cells.Item(32, 11).Value = 710.55054
$ws.Cells.Item(32, 12).Value = 1657.4546
$ws.Cells.Item(32, 13).Value = -423.55054
$ws.Cells.Item(32, 14).Value = -2231.4546

# Row 61 (ARM)
$ws.Cells.Item(61, 8).Value = 23303182
$ws.Cells.Item(61, 9).Value = 30333992
$ws.Cells.Item(61, 10).Value = 101513.8
$ws.Cells.Item(61, 11).Value = 30333992
$ws.Cells.Item(61, 12).Value = 101513.8
$ws.Cells.Item(61, 13).Value = -30333780
$ws.Cells.Item(61, 14).Value = -101937.8

# Row 96 (ARM)
$ws.Cells.Item(96, 8).Value = 28000
$ws.Cells.Item(96, 10).Value = 28000
$ws.Cells.Item(96, 12).Value = 28000
$ws.Cells.Item(96, 14).Value = -33492

# Row 132 (ARM)
$ws.Cells.Item(132, 8).Value = 78638.96000000001
$ws.Cells.Item(132, 9).Value = 50977.65
$ws.Cells.Item(132, 10).Value = 157671.28
$ws.Cells.Item(132, 11).Value = 152932.95
$ws.Cells.Item(132, 12).Value = 473013.84
$ws.Cells.Item(132, 13).Value = -150402.95
$ws.Cells.Item(132, 14).Value = -478073.84

# Row 136 (ARM)
$ws.Cells.Item(136, 8).Value = 23303182
$ws.Cells.Item(136, 9).Value = 30333992
$ws.Cells.Item(136, 10).Value = 101513.8
$ws.Cells.Item(136, 11).Value = 91001976
$ws.Cells.Item(136, 12).Value = 304541.4
$ws.Cells.Item(136, 13).Value = -90999426
$ws.Cells.Item(136, 14).Value = -309641.4

# Row 138 (ARM)
$ws.Cells.Item(138, 8).Value = 39499.668
$ws.Cells.Item(138, 10).Value = 39499.668
$ws.Cells.Item(138, 12).Value = 39499.668
$ws.Cells.Item(138, 14).Value = -49779.668

# Row 140 (ARM)
$ws.Cells.Item(140, 8).Value = 28602
$ws.Cells.Item(140, 10).Value = 28602
$ws.Cells.Item(140, 12).Value = 28602
$ws.Cells.Item(140, 14).Value = -38962

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Cells.Item(20, 8).Value = 1901.8
$ws.Cells.Item(20, 9).Value = 1700
$ws.Cells.Item(20, 10).Value = 2036.3334
$ws.Cells.Item(20, 11).Value = 1700
$ws.Cells.Item(20, 12).Value = 2036.3334
$ws.Cells.Item(20, 13).Value = -1453
$ws.Cells.Item(20, 14).Value = -2530.3334

# Row 99 (BSM)
$ws.Cells.Item(99, 8).Value = 1355.5555
$ws.Cells.Item(99, 9).Value = 1356.25
$ws.Cells.Item(99, 10).Value = 1350
$ws.Cells.Item(99, 11).Value = 1356.25
$ws.Cells.Item(99, 12).Value = 1350
$ws.Cells.Item(99, 13).Value = 141.75
$ws.Cells.Item(99, 14).Value = -4346

# Row 134 (BSM)
$ws.Cells.Item(134, 8).Value = 2316.2068
$ws.Cells.Item(134, 9).Value = 1349.4117
$ws.Cells.Item(134, 10).Value = 3685.8333
$ws.Cells.Item(134, 11).Value = 4048.2351
$ws.Cells.Item(134, 12).Value = 11057.4999
$ws.Cells.Item(134, 13).Value = -1513.2351
$ws.Cells.Item(134, 14).Value = -16127.4999

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 1762.7812
$ws.Cells.Item(31, 9).Value = 1047.4783
$ws.Cells.Item(31, 11).Value = 1047.4783
$ws.Cells.Item(31, 13).Value = -752.4783

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 1762.7812
$ws.Cells.Item(34, 9).Value = 1047.4783
$ws.Cells.Item(34, 11).Value = 1047.4783
$ws.Cells.Item(34, 13).Value = -845.4783

# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 23257734
$ws.Cells.Item(58, 9).Value = 30304788
$ws.Cells.Item(58, 10).Value = 2455.9
$ws.Cells.Item(58, 11).Value = 30304788
$ws.Cells.Item(58, 12).Value = 2455.9
$ws.Cells.Item(58, 13).Value = -30304585
$ws.Cells.Item(58, 14).Value = -2861.9

# Row 86 (CRP)
$ws.Cells.Item(86, 8).Value = 2492
$ws.Cells.Item(86, 9).Value = 2360.1177
$ws.Cells.Item(86, 10).Value = 2812.2856
$ws.Cells.Item(86, 11).Value = 2360.1177
$ws.Cells.Item(86, 12).Value = 2812.2856
$ws.Cells.Item(86, 13).Value = -1237.1177
$ws.Cells.Item(86, 14).Value = -5058.2856

# Row 89 (CRP)
$ws.Cells.Item(89, 8).Value = 2492
$ws.Cells.Item(89, 9).Value = 2360.1177
$ws.Cells.Item(89, 10).Value = 2812.2856
$ws.Cells.Item(89, 11).Value = 11800.5885
$ws.Cells.Item(89, 12).Value = 14061.428
$ws.Cells.Item(89, 13).Value = -6184.588499999998
$ws.Cells.Item(89, 14).Value = -25293.428

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 53250.64
$ws.Cells.Item(132, 9).Value = 36364.277
$ws.Cells.Item(132, 10).Value = 102221.1
$ws.Cells.Item(132, 11).Value = 109092.831
$ws.Cells.Item(132, 12).Value = 306663.3
$ws.Cells.Item(132, 13).Value = -106562.831
$ws.Cells.Item(132, 14).Value = -311723.3

# Row 134 (CRP)
$ws.Cells.Item(134, 8).Value = 30670.244
$ws.Cells.Item(134, 9).Value = 1474.1613
$ws.Cells.Item(134, 10).Value = 181516.67
$ws.Cells.Item(134, 11).Value = 4422.4839
$ws.Cells.Item(134, 12).Value = 544550.01
$ws.Cells.Item(134, 13).Value = -1887.4839
$ws.Cells.Item(134, 14).Value = -549620.01

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 23257734
$ws.Cells.Item(136, 9).Value = 30304788
$ws.Cells.Item(136, 10).Value = 2455.9
$ws.Cells.Item(136, 11).Value = 90914364
$ws.Cells.Item(136, 12).Value = 7367.700000000001
$ws.Cells.Item(136, 13).Value = -90911814
$ws.Cells.Item(136, 14).Value = -12467.7

$ws = $wb.Worksheets.Item("CUL")
# Row 124 (CUL)
$ws.Cells.Item(124, 8).Value = 962.63635

# Row 131 (CUL)
$ws.Cells.Item(131, 8).Value = 1031.2333
$ws.Cells.Item(131, 10).Value = 1216.8695
$ws.Cells.Item(131, 12).Value = 3650.6085
$ws.Cells.Item(131, 14).Value = -13730.6085

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 68404.664
$ws.Cells.Item(132, 9).Value = 44945.305
$ws.Cells.Item(132, 10).Value = 145485.42
$ws.Cells.Item(132, 11).Value = 134835.915
$ws.Cells.Item(132, 12).Value = 436456.26
$ws.Cells.Item(132, 13).Value = -132305.915
$ws.Cells.Item(132, 14).Value = -441516.26

# Row 141 (GSM)
$ws.Cells.Item(141, 8).Value = 36793.332
$ws.Cells.Item(141, 10).Value = 36793.332
$ws.Cells.Item(141, 12).Value = 36793.332
$ws.Cells.Item(141, 14).Value = -47153.332

$ws = $wb.Worksheets.Item("LTW")
# Row 100 (LTW)
$ws.Cells.Item(100, 8).Value = 1839.9
$ws.Cells.Item(100, 9).Value = 1599.75
$ws.Cells.Item(100, 11).Value = 1599.75
$ws.Cells.Item(100, 13).Value = -1058.75

# Row 132 (LTW)
$ws.Cells.Item(132, 8).Value = 20506.1
$ws.Cells.Item(132, 9).Value = 10256.5
$ws.Cells.Item(132, 10).Value = 49790.668
$ws.Cells.Item(132, 11).Value = 30769.5
$ws.Cells.Item(132, 12).Value = 149372.004
$ws.Cells.Item(132, 13).Value = -28239.5
$ws.Cells.Item(132, 14).Value = -154432.004

$ws = $wb.Worksheets.Item("WVR")
# Row 5 (WVR)
$ws.Cells.Item(5, 8).Value = 5000
$ws.Cells.Item(5, 10).Value = 5000
$ws.Cells.Item(5, 12).Value = 5000
$ws.Cells.Item(5, 14).Value = -5224

# Row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 65307.645
$ws.Cells.Item(132, 9).Value = 56211.168
$ws.Cells.Item(132, 10).Value = 92597.086
$ws.Cells.Item(132, 11).Value = 168633.504
$ws.Cells.Item(132, 12).Value = 277791.258
$ws.Cells.Item(132, 13).Value = -166103.504
$ws.Cells.Item(132, 14).Value = -282851.258

# Row 141 (WVR)
$ws.Cells.Item(141, 8).Value = 74977.78
$ws.Cells.Item(141, 10).Value = 74977.78
$ws.Cells.Item(141, 12).Value = 74977.78
$ws.Cells.Item(141, 14).Value = -85337.78
